# Add 2022-Q4 data: duplicate the existing "2021-Q4" sheet (same structure as
# every quarter sheet), rename the duplicate to "2022-Q4", position it right
# after "总计", overwrite its fund row with the new quarter's figures, and add
# a matching summary row on "总计".

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $text) {
    # Force a literal text value (not Excel's automatic "looks like a number"
    # coercion) while leaving the cell's style untouched.
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

# --- 1. Duplicate the "2021-Q4" sheet; the copy lands immediately before it,
#        which is exactly where the new quarter belongs. ------------------
$src = $wb.Worksheets.Item("2021-Q4")
$src.Copy($src)
$new = $wb.Worksheets.Item("2021-Q4 (2)")
$new.Name = "2022-Q4"

# --- 2. Overwrite the new sheet's single fund row with 2022-Q4 figures. ---
Set-TextValue $new.Range("B2") "015245"
$new.Range("C2").Value = "南华丰汇混合"
Set-TextValue $new.Range("D2") "0.11"
Set-TextValue $new.Range("E2") "84.24"
Set-TextValue $new.Range("F2") "1.09"
Set-TextValue $new.Range("G2") "0.0012"
# H2 (仓位排名) is unchanged from the copied sheet.

# --- 3. Update the "总计" roll-up sheet: push the existing rows down one
#        slot and add the new 2022-Q4 summary row on top. -----------------
$total = $wb.Worksheets.Item("总计")

# New row 4 = what used to be row 3 (2020-Q4), copying its style first.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2020-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.48

# Row 3 becomes what used to be row 2 (2021-Q4).
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.07000000000000001

# Row 2 becomes the brand-new 2022-Q4 summary.
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 1
$total.Range("D2").Value = 0

# --- 4. Restore the original active tab ("2020-Q4") -----------------------
# Copying a sheet makes the new copy active; put the selection back so only
# the intended data actually changed.
$wb.Worksheets.Item("2020-Q4").Activate()
